# Pooh Points: normal 20260214
# Updates the "Players" sheet with refreshed live-game stats (status clock
# advances from 15:41 to 11:31 - 2nd Half, box-score numbers increase) and
# swaps two pairs of players who were entered in the wrong rows by their
# backup/starter slot (Tyler Tanner <-> Rashaun Agee on The Oddities,
# Ali Dibba <-> Chandler Bing and Mike James <-> Pop Isaacs on Undrafted).
# Also refreshes the "OwnerTotals" rollup sheet to match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Per-row column updates for the Players sheet (row number -> column -> new value)
$rowData = @{
    2  = @{ G="11:31 - 2nd Half"; H=21; I=17; J=8; P=23; Q=6; R=10; U=4; V=4 }
    3  = @{ G="11:31 - 2nd Half"; H=18; I=20; O=2; P=26; Q=9; R=14; V=2 }
    4  = @{ G="11:31 - 2nd Half"; I=23; N=1; P=30; Q=8; R=14; T=11; U=2; V=3 }
    5  = @{ G="11:31 - 2nd Half"; H=5; I=5; K=4; O=4; P=20; Q=2; R=7; S=1; T=6 }
    6  = @{ D="Rashaun Agee"; E="TA&M"; G="11:31 - 2nd Half"; H=9; I=7; J=9; K=0; L=0; M=2; O=0; P=24; Q=3; R=11; T=0; U=1; V=2 }
    7  = @{ D="Tyler Tanner"; E="VAN"; G="11:31 - 2nd Half"; H=6; I=4; J=3; K=4; L=3; M=0; N=1; O=4; P=28; Q=2; R=9; T=2; U=0; V=0 }
    8  = @{ G="11:31 - 2nd Half"; H=6; J=8; K=2; P=19 }
    9  = @{ G="11:31 - 2nd Half"; H=-2; R=6 }
    10 = @{ G="11:31 - 2nd Half"; H=17; I=14; J=4; N=2; O=2; P=28; U=5; V=6 }
    11 = @{ D="Chandler Bing"; E="VAN"; G="11:31 - 2nd Half"; H=11; I=9; L=0; M=1; P=25; U=2 }
    12 = @{ D="Ali Dibba"; E="TA&M"; G="11:31 - 2nd Half"; H=10; I=8; J=6; L=2; M=0; N=1; P=18; U=1; V=2 }
    13 = @{ G="11:31 - 2nd Half"; H=8; I=2; K=1; M=1; N=1; O=4; P=12; Q=1; R=1 }
    14 = @{ G="11:31 - 2nd Half"; H=8; I=9; O=5; P=14; Q=4; R=5; S=1; T=2 }
    15 = @{ D="Pop Isaacs"; E="TA&M"; G="11:31 - 2nd Half"; H=5; I=5; J=1; K=2; L=1; O=1; P=21; Q=2; R=6; S=1; T=3 }
    16 = @{ D="Mike James"; E="VAN"; G="11:31 - 2nd Half"; H=3; I=0; J=2; K=2; O=4; P=16; Q=0; R=1; S=0; T=1 }
    17 = @{ G="11:31 - 2nd Half" }
    18 = @{ G="11:31 - 2nd Half" }
    19 = @{ G="11:31 - 2nd Half"; P=16 }
    20 = @{ G="11:31 - 2nd Half" }
}

foreach ($rowNum in $rowData.Keys) {
    $cols = $rowData[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws1.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}

# Refresh the OwnerTotals rollup sheet: "Boozers Losers" now leads with 21
# (was Hilton Heads at 19), so the two owners swap row order; The Oddities'
# and G-Flop's totals also climb as their starters' box scores update.
$ws2.Range("A2").Value = "Boozers Losers"
$ws2.Range("B2").Value = 21
$ws2.Range("C2").Value = 1

$ws2.Range("A3").Value = "Hilton Heads"
$ws2.Range("B3").Value = 19
$ws2.Range("C3").Value = 1

$ws2.Range("B4").Value = 15

$ws2.Range("B5").Value = 5
